$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = 50
$ws.Range("A3").Value = 32.99999999999983
$ws.Range("A4").Value = 37
$ws.Range("A5").Value = 28.39999999999964
$ws.Range("A6").Value = 26
$ws.Range("A7").Value = 78.59999999999854
$ws.Range("A8").Value = 44.19999999999891
$ws.Range("A9").Value = 93.39999999999964
